# systemd template support added
# Normalize machine "name" values to lower-case (ST1/ST2/ST3 -> st1/st2/st3)
# and refresh the sheet's cell formatting / selection to match the
# generated template output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lower-case the "name" column values (E2:E4) -------------------------
$ws.Range("E2").Value2 = "st1"
$ws.Range("E3").Value2 = "st2"
$ws.Range("E4").Value2 = "st3"

# --- Re-apply the base cell style across the used range -------------------
# This mirrors the formatting refresh that ships with the systemd template
# support change: every populated cell (A1:E4) is explicitly (re)stamped
# with the workbook's "Normal" style, producing a dedicated cell format
# record (distinct from the implicit default) while keeping the same
# locked/general/bottom-aligned appearance.
$used = $ws.Range("A1:E4")

# Base font keeps using the ANSI charset.
$used.Font.Charset = 1

$used.Style = "Normal"

# --- Update the active selection to the newly populated name column -------
$ws.Range("E2:E4").Select()
